$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0
    3  = 3
    4  = 4
    5  = 6
    6  = 7
    7  = 9
    8  = 12
    9  = 14
    10 = 16
    11 = 18
    12 = 20
    13 = 21
    14 = 24
    15 = 65
    16 = 89
    17 = 106
    18 = 181
    19 = 190
    20 = 230
    21 = 253
    22 = 270
    23 = 329
    24 = 361
    25 = 379
    26 = 391
    27 = 447
}

foreach ($row in $values.Keys) {
    $ws.Range("A$row").Value = $values[$row]
}
